$d = $word.ActiveDocument

# Remove the placeholder text "vnpt.SiteAddress" (and its run) that followed
# "Địa chỉ: " in the first occurrence of that label.
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 2)
